# New crime data collected - weekly CompStat update for the 47th Precinct
# Updates: report week/volume header text, and the weekly/28-day/YTD/2yr crime
# statistics table (rows 14-29) with newly collected figures.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# ---------------------------------------------------------------------------
# Helper: replace a single occurrence of $oldSub with $newSub inside a cell's
# text, preserving the surrounding rich-text run (used for the "Volume NN
# Number NN" and "Report Covering the Week ... Through ..." header cells,
# which are shared strings made up of multiple runs).
# ---------------------------------------------------------------------------
function Set-SubstringInCell {
    param($range, [string]$oldSub, [string]$newSub)
    $full = $range.Value2
    $idx = $full.IndexOf($oldSub)
    if ($idx -lt 0) {
        throw "Substring '$oldSub' not found in '$full'"
    }
    $chars = $range.Characters($idx + 1, $oldSub.Length)
    $chars.Text = $newSub
}

# ---------------------------------------------------------------------------
# Helper: turn a cell into a text "placeholder" cell (e.g. "0" or "***.*")
# by copying both the number format and the literal text value from a
# reference cell that already uses that placeholder (format + shared-string
# text together, so the destination becomes a genuine text cell instead of
# Excel auto-converting "0" back into a number). This mirrors how blank/NA
# data points are represented elsewhere in this report.
# ---------------------------------------------------------------------------
function Set-TextPlaceholder {
    param($dst, $src)
    $src.Copy() | Out-Null
    $dst.PasteSpecial(-4122) | Out-Null  # xlPasteFormats
    $src.Copy() | Out-Null
    $dst.PasteSpecial(-4163) | Out-Null  # xlPasteValues
    $excel.CutCopyMode = 0
}

# ---------------------------------------------------------------------------
# Header: Volume / Number and the reporting week range
# ---------------------------------------------------------------------------
Set-SubstringInCell $ws.Range("A8") "39" "40"
Set-SubstringInCell $ws.Range("C9") "9/25/2023" "10/2/2023"
Set-SubstringInCell $ws.Range("C9") "10/1/2023" "10/8/2023"

# ---------------------------------------------------------------------------
# Row 14 - Murder
# ---------------------------------------------------------------------------
$ws.Range("N14").Value = -58.333333333333

# ---------------------------------------------------------------------------
# Row 15 - Rape
# ---------------------------------------------------------------------------
$ws.Range("I15").Value = 33
$ws.Range("K15").Value = -15.384615384615
$ws.Range("L15").Value = 3.125
$ws.Range("M15").Value = 3.125
$ws.Range("N15").Value = -45

# ---------------------------------------------------------------------------
# Row 16 - Robbery
# ---------------------------------------------------------------------------
$ws.Range("C16").Value = 7
$ws.Range("D16").Value = 5
$ws.Range("E16").Value = 40
$ws.Range("F16").Value = 32
$ws.Range("G16").Value = 32
$ws.Range("H16").Value = 0
$ws.Range("I16").Value = 359
$ws.Range("J16").Value = 341
$ws.Range("K16").Value = 5.278592375366
$ws.Range("L16").Value = 34.456928838951
$ws.Range("M16").Value = 10.802469135802
$ws.Range("N16").Value = -67.033976124885

# ---------------------------------------------------------------------------
# Row 17 - Fel. Assault
# ---------------------------------------------------------------------------
$ws.Range("C17").Value = 21
$ws.Range("D17").Value = 12
$ws.Range("E17").Value = 75
$ws.Range("F17").Value = 71
$ws.Range("G17").Value = 49
$ws.Range("H17").Value = 44.897959183673
$ws.Range("I17").Value = 656
$ws.Range("J17").Value = 574
$ws.Range("K17").Value = 14.285714285714
$ws.Range("L17").Value = 17.142857142857
$ws.Range("M17").Value = 99.392097264437
$ws.Range("N17").Value = -5.065123010130

# ---------------------------------------------------------------------------
# Row 18 - Burglary
# ---------------------------------------------------------------------------
$ws.Range("C18").Value = 2
$ws.Range("D18").Value = 10
$ws.Range("E18").Value = -80
$ws.Range("G18").Value = 22
$ws.Range("H18").Value = -40.909090909090
$ws.Range("I18").Value = 208
$ws.Range("J18").Value = 224
$ws.Range("K18").Value = -7.142857142857
$ws.Range("L18").Value = 11.229946524064
$ws.Range("M18").Value = -22.388059701492
$ws.Range("N18").Value = -85.753424657534

# ---------------------------------------------------------------------------
# Row 19 - Gr. Larceny
# ---------------------------------------------------------------------------
$ws.Range("C19").Value = 21
$ws.Range("D19").Value = 21
$ws.Range("E19").Value = 0
$ws.Range("F19").Value = 70
$ws.Range("G19").Value = 71
$ws.Range("H19").Value = -1.408450704225
$ws.Range("I19").Value = 594
$ws.Range("J19").Value = 592
$ws.Range("K19").Value = 0.337837837837
$ws.Range("L19").Value = 51.530612244898
$ws.Range("M19").Value = 138.55421686747
$ws.Range("N19").Value = 37.819025522041

# ---------------------------------------------------------------------------
# Row 20 - G.L.A.
# ---------------------------------------------------------------------------
$ws.Range("C20").Value = 13
$ws.Range("D20").Value = 8
$ws.Range("E20").Value = 62.5
$ws.Range("F20").Value = 61
$ws.Range("G20").Value = 26
$ws.Range("H20").Value = 134.615384615385
$ws.Range("I20").Value = 499
$ws.Range("J20").Value = 337
$ws.Range("K20").Value = 48.071216617210
$ws.Range("L20").Value = 29.610389610389
$ws.Range("M20").Value = 98.015873015873
$ws.Range("N20").Value = -60.239043824701

# ---------------------------------------------------------------------------
# Row 21 - TOTAL
# ---------------------------------------------------------------------------
$ws.Range("C21").Value = 65
$ws.Range("D21").Value = 56
$ws.Range("E21").Value = 16.071428571428
$ws.Range("F21").Value = 249
$ws.Range("G21").Value = 202
$ws.Range("H21").Value = 23.267326732673
$ws.Range("I21").Value = 2359
$ws.Range("J21").Value = 2118
$ws.Range("K21").Value = 11.378659112370
$ws.Range("L21").Value = 28.485838779956
$ws.Range("M21").Value = 60.149355057705
$ws.Range("N21").Value = -52.914171656686

# ---------------------------------------------------------------------------
# Row 22 - Transit
# ---------------------------------------------------------------------------
$ws.Range("C22").Value = 2
$ws.Range("E22").Value = 100
$ws.Range("F22").Value = 6
$ws.Range("G22").Value = 3
$ws.Range("I22").Value = 22
$ws.Range("J22").Value = 33
$ws.Range("K22").Value = -33.333333333333
$ws.Range("L22").Value = 100
$ws.Range("M22").Value = 10

# ---------------------------------------------------------------------------
# Row 23 - Housing (C23 becomes a real number this week instead of "0")
# ---------------------------------------------------------------------------
$ws.Range("D22").Copy() | Out-Null
$ws.Range("C23").PasteSpecial(-4122) | Out-Null  # xlPasteFormats (numeric style)
$excel.CutCopyMode = 0
$ws.Range("C23").Value = 1
$ws.Range("D23").Value = 1
$ws.Range("E23").Value = 0
$ws.Range("F23").Value = 5
$ws.Range("G23").Value = 8
$ws.Range("H23").Value = -37.5
$ws.Range("I23").Value = 81
$ws.Range("J23").Value = 90
$ws.Range("K23").Value = -10
$ws.Range("L23").Value = -6.896551724137
$ws.Range("M23").Value = 32.786885245901

# ---------------------------------------------------------------------------
# Row 24 - Petit Larceny
# ---------------------------------------------------------------------------
$ws.Range("C24").Value = 19
$ws.Range("D24").Value = 40
$ws.Range("E24").Value = -52.5
$ws.Range("F24").Value = 96
$ws.Range("G24").Value = 118
$ws.Range("H24").Value = -18.644067796610
$ws.Range("I24").Value = 1034
$ws.Range("J24").Value = 1173
$ws.Range("K24").Value = -11.849957374254
$ws.Range("L24").Value = 23.536439665471
$ws.Range("M24").Value = 78.275862068965

# ---------------------------------------------------------------------------
# Row 25 - Misd. Assault
# ---------------------------------------------------------------------------
$ws.Range("C25").Value = 24
$ws.Range("D25").Value = 17
$ws.Range("E25").Value = 41.176470588235
$ws.Range("F25").Value = 77
$ws.Range("G25").Value = 69
$ws.Range("H25").Value = 11.594202898550
$ws.Range("I25").Value = 739
$ws.Range("J25").Value = 699
$ws.Range("K25").Value = 5.722460658082
$ws.Range("L25").Value = 3.646563814866
$ws.Range("M25").Value = -3.145478374836

# ---------------------------------------------------------------------------
# Row 26 - UCR Rape*
# ---------------------------------------------------------------------------
$ws.Range("C26").Value = 3
$ws.Range("F26").Value = 11
$ws.Range("H26").Value = 83.333333333333
$ws.Range("I26").Value = 55
$ws.Range("K26").Value = -15.384615384615
$ws.Range("L26").Value = 14.583333333333

# ---------------------------------------------------------------------------
# Row 27 - Other Sex Crimes
# ---------------------------------------------------------------------------
$ws.Range("C27").Value = 1
$ws.Range("D27").Value = 1
$ws.Range("E27").Value = 0
$ws.Range("I27").Value = 63
$ws.Range("J27").Value = 60
$ws.Range("K27").Value = 5
$ws.Range("L27").Value = 5

# ---------------------------------------------------------------------------
# Row 28 - Shooting Vic. (D28/E28 become "0"/"***.*" placeholders this week)
# ---------------------------------------------------------------------------
$placeholderZero = $ws.Range("C28")
$placeholderPct = $ws.Range("E14")
Set-TextPlaceholder $ws.Range("D28") $placeholderZero
Set-TextPlaceholder $ws.Range("E28") $placeholderPct
$ws.Range("F28").Value = 2
$ws.Range("H28").Value = 0
$ws.Range("L28").Value = -41.379310344827

# ---------------------------------------------------------------------------
# Row 29 - Shooting Inc. (D29/E29 become "0"/"***.*" placeholders this week)
# ---------------------------------------------------------------------------
Set-TextPlaceholder $ws.Range("D29") $placeholderZero
Set-TextPlaceholder $ws.Range("E29") $placeholderPct
$ws.Range("F29").Value = 2
$ws.Range("H29").Value = 0
$ws.Range("L29").Value = -40.816326530612
